$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.055.96'
$ws.Range("E2").Value = '  +4.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.228.27'
$ws.Range("E3").Value = '  +2.99%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.13'
$ws.Range("E5").Value = '  +2.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '82.75'
$ws.Range("E6").Value = '  +13.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +3.21%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.611'
$ws.Range("E9").Value = '  +5.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.15'
$ws.Range("E10").Value = '  +11.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0933'
$ws.Range("E11").Value = '  +2.82%  '
$ws.Range("E12").Value = '  +5.17%  '
$ws.Range("E13").Value = '  +2.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.565.59'
$ws.Range("E14").Value = '  +3.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.68'
$ws.Range("E15").Value = '  +3.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.237.43'
$ws.Range("E16").Value = '  +4.46%  '
$ws.Range("E17").Value = '  +3.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.923.20'
$ws.Range("E18").Value = '  +4.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000104'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.45'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("E21").Value = '  +4.02%  '
$ws.Range("E22").Value = '  +11.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.19'
$ws.Range("E23").Value = '  +3.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.25'
$ws.Range("E24").Value = '  -2.89%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  +3.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '40.68'
$ws.Range("E27").Value = '  +11.80%  '
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +2.65%  '
$ws.Range("E30").Value = '  +4.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.14'
$ws.Range("E31").Value = '  +2.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0902'
$ws.Range("E32").Value = '  +12.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.67'
$ws.Range("E33").Value = '  +3.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.35'
$ws.Range("E34").Value = '  +4.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.117'
$ws.Range("E35").Value = '  +9.56%  '
$ws.Range("E36").Value = '  +2.63%  '
$ws.Range("E37").Value = '  +11.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.52'
$ws.Range("E38").Value = '  +7.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.43'
$ws.Range("E39").Value = '  +13.99%  '
$ws.Range("E40").Value = '  +26.18%  '
$ws.Range("E41").Value = '  +4.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.04'
$ws.Range("E42").Value = '  +9.49%  '
$ws.Range("E43").Value = '  +8.15%  '
$ws.Range("E44").Value = '  +4.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.67'
$ws.Range("E45").Value = '  +1.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0989'
$ws.Range("E46").Value = '  +2.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.37'
$ws.Range("E47").Value = '  +2.31%  '
$ws.Range("E48").Value = '  +30.14%  '
$ws.Range("E49").Value = '  +4.56%  '
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("E51").Value = '  +3.34%  '
